# Apply the Entsoe Actual Production Hydro Water Reservoir data refresh:
# the series rolls forward by one day (row 2 now starts at what was
# row 98, i.e. 2025-09-xx+1) and a new day of 15-minute readings is
# appended at the end (rows 98:193).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: timestamps (Excel serial dates), column B: actual production (MW)
$timestamps = @(
    45918, 45918.01041666666, 45918.02083333334, 45918.03125, 45918.04166666666, 45918.05208333334, 45918.0625, 45918.07291666666, 45918.08333333334, 45918.09375, 45918.10416666666, 45918.11458333334, 45918.125, 45918.13541666666, 45918.14583333334, 45918.15625, 45918.16666666666, 45918.17708333334, 45918.1875, 45918.19791666666, 45918.20833333334, 45918.21875, 45918.22916666666, 45918.23958333334, 45918.25, 45918.26041666666, 45918.27083333334, 45918.28125, 45918.29166666666, 45918.30208333334, 45918.3125, 45918.32291666666, 45918.33333333334, 45918.34375, 45918.35416666666, 45918.36458333334, 45918.375, 45918.38541666666, 45918.39583333334, 45918.40625, 45918.41666666666, 45918.42708333334, 45918.4375, 45918.44791666666, 45918.45833333334, 45918.46875, 45918.47916666666, 45918.48958333334, 45918.5, 45918.51041666666, 45918.52083333334, 45918.53125, 45918.54166666666, 45918.55208333334, 45918.5625, 45918.57291666666, 45918.58333333334, 45918.59375, 45918.60416666666, 45918.61458333334, 45918.625, 45918.63541666666, 45918.64583333334, 45918.65625, 45918.66666666666, 45918.67708333334, 45918.6875, 45918.69791666666, 45918.70833333334, 45918.71875, 45918.72916666666, 45918.73958333334, 45918.75, 45918.76041666666, 45918.77083333334, 45918.78125, 45918.79166666666, 45918.80208333334, 45918.8125, 45918.82291666666, 45918.83333333334, 45918.84375, 45918.85416666666, 45918.86458333334, 45918.875, 45918.88541666666, 45918.89583333334, 45918.90625, 45918.91666666666, 45918.92708333334, 45918.9375, 45918.94791666666, 45918.95833333334, 45918.96875, 45918.97916666666, 45918.98958333334, 45919, 45919.01041666666, 45919.02083333334, 45919.03125, 45919.04166666666, 45919.05208333334, 45919.0625, 45919.07291666666, 45919.08333333334, 45919.09375, 45919.10416666666, 45919.11458333334, 45919.125, 45919.13541666666, 45919.14583333334, 45919.15625, 45919.16666666666, 45919.17708333334, 45919.1875, 45919.19791666666, 45919.20833333334, 45919.21875, 45919.22916666666, 45919.23958333334, 45919.25, 45919.26041666666, 45919.27083333334, 45919.28125, 45919.29166666666, 45919.30208333334, 45919.3125, 45919.32291666666, 45919.33333333334, 45919.34375, 45919.35416666666, 45919.36458333334, 45919.375, 45919.38541666666, 45919.39583333334, 45919.40625, 45919.41666666666, 45919.42708333334, 45919.4375, 45919.44791666666, 45919.45833333334, 45919.46875, 45919.47916666666, 45919.48958333334, 45919.5, 45919.51041666666, 45919.52083333334, 45919.53125, 45919.54166666666, 45919.55208333334, 45919.5625, 45919.57291666666, 45919.58333333334, 45919.59375, 45919.60416666666, 45919.61458333334, 45919.625, 45919.63541666666, 45919.64583333334, 45919.65625, 45919.66666666666, 45919.67708333334, 45919.6875, 45919.69791666666, 45919.70833333334, 45919.71875, 45919.72916666666, 45919.73958333334, 45919.75, 45919.76041666666, 45919.77083333334, 45919.78125, 45919.79166666666, 45919.80208333334, 45919.8125, 45919.82291666666, 45919.83333333334, 45919.84375, 45919.85416666666, 45919.86458333334, 45919.875, 45919.88541666666, 45919.89583333334, 45919.90625, 45919.91666666666, 45919.92708333334, 45919.9375, 45919.94791666666, 45919.95833333334, 45919.96875, 45919.97916666666, 45919.98958333334
)
$production = @(
    365, 362, 360, 356, 360, 358, 359, 358, 356, 351, 353, 351, 353, 351, 352, 351, 356, 356, 355, 361, 460, 463, 463, 467, 505, 501, 502, 503, 411, 391, 390, 373, 399, 373, 371, 366, 150, 100, 84, 74, 50, 39, 38, 32, 93, 89, 83, 84, 66, 80, 86, 82, 45, 66, 66, 65, 66, 67, 68, 67, 84, 89, 90, 105, 464, 492, 499, 523, 621, 650, 652, 664, 810, 819, 819, 836, 876, 885, 882, 881, 869, 864, 860, 863, 776, 769, 770, 762, 462, 446, 448, 444, 396, 389, 388, 385, 350, 343, 338, 335, 342, 337, 339, 335, 355, 355, 356, 355, 355, 355, 353, 354, 376, 385, 387, 393, 471, 475, 478, 477, 509, 510, 509, 507, 536, 523, 510, 501, 503, 468, 470, 434, 366, 358, 354, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0
)

$rowCount = 192
$data = New-Object 'object[,]' $rowCount,2
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = $timestamps[$i]
    $data[$i,1] = $production[$i]
}

$ws.Range("A2:B193").Value = $data
